# Generate Report for Handback
# The c617b04f-6424-4b11-adb3-db1762f1a537.md file has now been handed back
# (in sync with en-US) for both the zh-cn and de-de locales. Update the
# status, handback timestamps and clear the stale "version mismatch" error
# on the per-locale sheets, and reflect the new status on the Overview sheet.

$wb = $excel.ActiveWorkbook

# ----- Overview sheet -----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# ----- zh-cn sheet -----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("K3").Value = "2016-09-05 18:55:13"
$zhcn.Range("P3").Value = ""
$zhcn.Columns.Item(16).ColumnWidth = 12.83

# ----- de-de sheet -----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("K3").Value = "2016-09-05 18:55:21"
$dede.Range("P3").Value = ""
$dede.Columns.Item(16).ColumnWidth = 12.83
